$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.637.52"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.786.28"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'352.56"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").Value = "'109.10"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  -2.40%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("D10").Value = "'39.65"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("D13").Value = "'20.06"
$ws.Range("E13").Value = "  +2.74%  "
$ws.Range("D14").Value = "'7.68"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").Value = "3.222.09"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").Value = "2.795.92"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("D17").Value = "'0.928"
$ws.Range("E17").Value = "  -2.66%  "
$ws.Range("D18").Value = "51.625.30"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "'7.71"
$ws.Range("E19").Value = "  +3.76%  "
$ws.Range("D20").Value = "'3.18"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").Value = "'267.32"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").Value = "'2.74"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").Value = "'26.10"
$ws.Range("E26").Value = "  -2.42%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  +12.81%  "
$ws.Range("D29").Value = "'10.26"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").Value = "'37.02"
$ws.Range("E30").Value = "  +7.32%  "
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("E32").Value = "  +6.85%  "
$ws.Range("D33").Value = "'51.74"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").Value = "'0.0453"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "'5.67"
$ws.Range("E35").Value = "  +7.29%  "
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "'18.53"
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("E39").Value = "  -3.21%  "
$ws.Range("E40").Value = "  -1.97%  "
$ws.Range("E41").Value = "  -1.86%  "
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "'120.40"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").Value = "'22.08"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("D46").Value = "2.124.67"
$ws.Range("E46").Value = "  +1.98%  "
$ws.Range("D47").Value = "'3.34"
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("E48").Value = "  +4.47%  "
$ws.Range("B49").Value = "SEI"
$ws.Range("C49").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D49").Value = "'0.907"
$ws.Range("E49").Value = "  -3.01%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'5.43"
$ws.Range("E50").Value = "  -5.26%  "
$ws.Range("E51").Value = "  +8.27%  "
